$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notenrechner")

# Fill in the points for the reviewed "FotografInnen" (photographers) criteria
$ws.Range("B11").Value = 1
$ws.Range("B14").Value = 1
$ws.Range("B15").Value = 2
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 1

# Move the active selection to reflect where the edit was made
$ws.Activate()
$ws.Range("A17").Select()
